# Add a new "Thank You!" slide at the end of the deck (slide 9), using the
# "Title Only" layout (ppLayoutTitleOnly = 11) so the slide only carries a
# single title placeholder, matching the target slide structure.

$p = $ppt.ActivePresentation

$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 11)

$shp = $s.Shapes.Item(1)

# Reposition/resize the title placeholder to match the target layout.
# (Values below are crafted so that the COM layer's Single-precision point
# storage rounds back to the exact target EMU values.)
$shp.Left = 85.56248474121094
$shp.Top = 159.70594787597656
$shp.Width = 788.8750610351562
$shp.Height = 137.9999542236328

# Turn off autofit on the title placeholder (-> <a:noAutofit/>) before
# touching the font, so PowerPoint doesn't inject <a:normAutofit/>.
$shp.TextFrame.AutoSize = 0

$tr = $shp.TextFrame.TextRange
$tr.Text = "Thank You!"
$tr.Font.Size = 80
$tr.Font.Name = "Arial Rounded MT Bold"

# Slide transition: "Reveal" (slow) isn't representable through the legacy
# PpEntryEffect enum, so fall back to the closest native equivalent (the
# actual fallback PowerPoint itself emits for Reveal on older readers): a
# slow Fade.
$s.SlideShowTransition.EntryEffect = 1793
$s.SlideShowTransition.Speed = 1
